$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.6606524410359556; C = 1.655778082260271;  D = 3.537761648806719;  E = 0.4942365360607697;  G = 6.348428708163715 }
    3 = @{ B = 0.1190320826869504; C = 10.34677158129881;   D = 3.537761648806719;  E = 10.19245300693656;   G = 24.19601831972904 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271;   D = 3.537761648806719;  E = 0.4942365360607697;  G = 8.974608811992548 }
    5 = @{ B = 3.286832544864788;  C = 1.655778082260271;   D = 0.7527432677738641; E = 0.4942365360607697;  G = 6.189590430959694 }
    6 = @{ B = 1.455362044514542;  C = 0.306821227259698;   D = 0.1494219747398047; E = 0.4942365360607697;  G = 2.405841782574814 }
    7 = @{ B = 3.286832544864788;  C = 1.655778082260271;   D = 0.1494219747398047; E = 0.4942365360607697;  G = 5.586269137925634 }
    8 = @{ B = 3.286832544864788;  C = 1.655778082260271;   D = 0.7527432677738641; E = 0.4942365360607697;  G = 6.189590430959694 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
